# Scheduled runner update: refresh computed market-board profit figures
# (columns H-N) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3050.7727
$ws.Range("I64").Value = 3051
$ws.Range("J64").Value = 3050.6428
$ws.Range("K64").Value = 3051
$ws.Range("L64").Value = 3050.6428
$ws.Range("M64").Value = -2803
$ws.Range("N64").Value = -3546.6428
$ws.Range("H67").Value = 3050.7727
$ws.Range("I67").Value = 3051
$ws.Range("J67").Value = 3050.6428
$ws.Range("K67").Value = 3051
$ws.Range("L67").Value = 3050.6428
$ws.Range("M67").Value = -2193
$ws.Range("N67").Value = -4766.6428
$ws.Range("H70").Value = 1307.8077
$ws.Range("I70").Value = 1583.4
$ws.Range("J70").Value = 932
$ws.Range("K70").Value = 4750.200000000001
$ws.Range("L70").Value = 2796
$ws.Range("M70").Value = -4480.200000000001
$ws.Range("N70").Value = -3336
$ws.Range("H73").Value = 1307.8077
$ws.Range("I73").Value = 1583.4
$ws.Range("J73").Value = 932
$ws.Range("K73").Value = 4750.200000000001
$ws.Range("L73").Value = 2796
$ws.Range("M73").Value = -3814.200000000001
$ws.Range("N73").Value = -4668
$ws.Range("H76").Value = 373952.5
$ws.Range("J76").Value = 4298.6
$ws.Range("L76").Value = 4298.6
$ws.Range("N76").Value = -4928.6
$ws.Range("H79").Value = 373952.5
$ws.Range("J79").Value = 4298.6
$ws.Range("L79").Value = 4298.6
$ws.Range("N79").Value = -6482.6
$ws.Range("H138").Value = 3055.4285
$ws.Range("I138").Value = 1671.9048
$ws.Range("J138").Value = 3885.543
$ws.Range("K138").Value = 5015.7144
$ws.Range("L138").Value = 11656.629
$ws.Range("M138").Value = 124.2856000000002
$ws.Range("N138").Value = -21936.629

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 948.2414
$ws.Range("I2").Value = 832.26666
$ws.Range("J2").Value = 1072.5
$ws.Range("K2").Value = 832.26666
$ws.Range("L2").Value = 1072.5
$ws.Range("M2").Value = -719.26666
$ws.Range("N2").Value = -1298.5
$ws.Range("H110").Value = 2761
$ws.Range("J110").Value = 2685
$ws.Range("L110").Value = 2685
$ws.Range("N110").Value = -6775
$ws.Range("H116").Value = 948.2414
$ws.Range("I116").Value = 832.26666
$ws.Range("J116").Value = 1072.5
$ws.Range("K116").Value = 832.26666
$ws.Range("L116").Value = 1072.5
$ws.Range("M116").Value = 1461.73334
$ws.Range("N116").Value = -5660.5
$ws.Range("H122").Value = 3233.4736
$ws.Range("I122").Value = 3048.4
$ws.Range("J122").Value = 3439.111
$ws.Range("K122").Value = 9145.200000000001
$ws.Range("L122").Value = 10317.333
$ws.Range("M122").Value = -6695.200000000001
$ws.Range("N122").Value = -15217.333
$ws.Range("H132").Value = 2524.9678
$ws.Range("I132").Value = 1994.7826
$ws.Range("J132").Value = 4049.25
$ws.Range("K132").Value = 5984.3478
$ws.Range("L132").Value = 12147.75
$ws.Range("M132").Value = -3454.3478
$ws.Range("N132").Value = -17207.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 948.2414
$ws.Range("I3").Value = 832.26666
$ws.Range("J3").Value = 1072.5
$ws.Range("K3").Value = 832.26666
$ws.Range("L3").Value = 1072.5
$ws.Range("M3").Value = -718.26666
$ws.Range("N3").Value = -1300.5
$ws.Range("H20").Value = 931.3913
$ws.Range("I20").Value = 690.8570999999999
$ws.Range("J20").Value = 1305.5555
$ws.Range("K20").Value = 690.8570999999999
$ws.Range("L20").Value = 1305.5555
$ws.Range("M20").Value = -443.8570999999999
$ws.Range("N20").Value = -1799.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 419087.25
$ws.Range("I86").Value = 835482.8
$ws.Range("J86").Value = 2691.6667
$ws.Range("K86").Value = 835482.8
$ws.Range("L86").Value = 2691.6667
$ws.Range("M86").Value = -834359.8
$ws.Range("N86").Value = -4937.6667
$ws.Range("H89").Value = 419087.25
$ws.Range("I89").Value = 835482.8
$ws.Range("J89").Value = 2691.6667
$ws.Range("K89").Value = 4177414
$ws.Range("L89").Value = 13458.3335
$ws.Range("M89").Value = -4171798
$ws.Range("N89").Value = -24690.3335
$ws.Range("H132").Value = 62502468
$ws.Range("I132").Value = 166670050
$ws.Range("J132").Value = 1914.6
$ws.Range("K132").Value = 500010150
$ws.Range("L132").Value = 5743.799999999999
$ws.Range("M132").Value = -500007620
$ws.Range("N132").Value = -10803.8
$ws.Range("H134").Value = 3719.3225
$ws.Range("I134").Value = 3592.2856
$ws.Range("J134").Value = 4905
$ws.Range("K134").Value = 10776.8568
$ws.Range("L134").Value = 14715
$ws.Range("M134").Value = -8241.856800000001
$ws.Range("N134").Value = -19785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 43657140
$ws.Range("I70").Value = 118490280
$ws.Range("J70").Value = 4475
$ws.Range("K70").Value = 118490280
$ws.Range("L70").Value = 4475
$ws.Range("M70").Value = -118490010
$ws.Range("N70").Value = -5015
$ws.Range("H73").Value = 43657140
$ws.Range("I73").Value = 118490280
$ws.Range("J73").Value = 4475
$ws.Range("K73").Value = 118490280
$ws.Range("L73").Value = 4475
$ws.Range("M73").Value = -118489344
$ws.Range("N73").Value = -6347
$ws.Range("H80").Value = 2421.3635
$ws.Range("I80").Value = 2015
$ws.Range("J80").Value = 2611
$ws.Range("K80").Value = 2015
$ws.Range("L80").Value = 2611
$ws.Range("M80").Value = -1017
$ws.Range("N80").Value = -4607
$ws.Range("H83").Value = 2421.3635
$ws.Range("I83").Value = 2015
$ws.Range("J83").Value = 2611
$ws.Range("K83").Value = 10075
$ws.Range("L83").Value = 13055
$ws.Range("M83").Value = -5083
$ws.Range("N83").Value = -23039
$ws.Range("H102").Value = 1689
$ws.Range("I102").Value = 1800.7778
$ws.Range("J102").Value = 1521.3334
$ws.Range("K102").Value = 1800.7778
$ws.Range("L102").Value = 1521.3334
$ws.Range("M102").Value = -178.7778000000001
$ws.Range("N102").Value = -4765.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1648.4193
$ws.Range("I7").Value = 1315.7778
$ws.Range("J7").Value = 1784.5
$ws.Range("K7").Value = 1315.7778
$ws.Range("L7").Value = 1784.5
$ws.Range("M7").Value = -1203.7778
$ws.Range("N7").Value = -2008.5
$ws.Range("H126").Value = 1648.4193
$ws.Range("I126").Value = 1315.7778
$ws.Range("J126").Value = 1784.5
$ws.Range("K126").Value = 3947.3334
$ws.Range("L126").Value = 5353.5
$ws.Range("M126").Value = -1477.3334
$ws.Range("N126").Value = -10293.5
$ws.Range("H132").Value = 8777583
$ws.Range("I132").Value = 13520536
$ws.Range("J132").Value = 3121.4
$ws.Range("K132").Value = 40561608
$ws.Range("L132").Value = 9364.200000000001
$ws.Range("M132").Value = -40559078
$ws.Range("N132").Value = -14424.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2333.5715
$ws.Range("I126").Value = 2788
$ws.Range("K126").Value = 8364
$ws.Range("M126").Value = -5894
$ws.Range("H132").Value = 2381.261
$ws.Range("I132").Value = 1015.1
$ws.Range("J132").Value = 3432.1538
$ws.Range("K132").Value = 3045.3
$ws.Range("L132").Value = 10296.4614
$ws.Range("M132").Value = -515.3000000000002
$ws.Range("N132").Value = -15356.4614

